$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update D2 (poucent / "combien de dollars j'ajoute a chaque achat") from 1 to 2 ---
$ws.Range("D2").Value = 2

# --- Row 16: give it the same row style/height as the rows above it (it is no longer the newest row) ---
$ws.Range("A16:D16").ClearContents()
$ws.Range("A15:D15").Copy()
$ws.Range("A16:D16").Select() | Out-Null
$ws.Paste() | Out-Null
$ws.Range("A16").Value = 0.9431
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "19/11/2025"
$ws.Range("D16").Value = "19:01:43"
$ws.Rows.Item(16).RowHeight = 12.75

# --- Row 17: new DCA entry, styled like the rows above ---
$ws.Range("A17:D17").ClearContents()
$ws.Range("A15:D15").Copy()
$ws.Range("A17:D17").Select() | Out-Null
$ws.Paste() | Out-Null
$ws.Range("A17").Value = 0.914
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = "20/11/2025"
$ws.Range("D17").Value = "19:01:20"
$ws.Rows.Item(17).RowHeight = 12.75

# --- Row 18: new DCA entry, styled like the rows above ---
$ws.Range("A18:D18").ClearContents()
$ws.Range("A15:D15").Copy()
$ws.Range("A18:D18").Select() | Out-Null
$ws.Paste() | Out-Null
$ws.Range("A18").Value = 0.9115
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "20/11/2025"
$ws.Range("D18").Value = "20:14:55"
$ws.Rows.Item(18).RowHeight = 12.75

# --- Row 19: newest DCA entry, kept unstyled (default) like row 16 used to be ---
$ws.Range("A19").Value = 0.8855
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = "21/11/2025"
$ws.Range("D19").Value = "14:01:31"

$excel.CutCopyMode = $false

# --- Move the live selection to D2 ---
$ws.Range("D2").Select() | Out-Null

Write-Output "edit applied"
